# Apply the "storeKeys" / column-realignment edit to the '#system' sheet
# (sheet1, the hidden lookup table that backs the nexial expression
# autocomplete named ranges).
#
# Summary of the change:
#  1. A new JSON expression "storeKeys(json,jsonpath,var)" is inserted
#     (alphabetically) into the "json" function list in column M, between
#     "storeCount(...)" (M15) and "storeValue(...)" (old M16). Everything
#     from old M16 downward shifts down by one row.
#  2. The single-cell "text" column (column Y, header "text") is removed.
#     Columns Z..AE shift left by one column (Z->Y, AA->Z, AB->AA, AC->AB,
#     AD->AC, AE->AD) and the vacated column AE is cleared.
#  3. The "text" category name itself disappears from the category list in
#     column A (used to be at A25); rows A26:A31 shift up to A25:A30 and
#     A31 is cleared.
#  4. The workbook-level defined names that describe these ranges are
#     updated to match the new extents. (The separate "text" named range
#     itself intentionally keeps pointing at $Y$2:$Y$2 -- unrelated to the
#     "text" category that disappeared from column A -- matching upstream.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Insert "storeKeys(json,jsonpath,var)" into the json list (column M).
#    Shift M16:M17 down to M17:M18, bottom-up so nothing gets clobbered.
# ---------------------------------------------------------------------
for ($row = 18; $row -ge 17; $row--) {
    $srcRow = $row - 1
    $val = $ws.Cells.Item($srcRow, 13).Value()
    $ws.Cells.Item($row, 13).Value = $val
}
$ws.Cells.Item(16, 13).Value = "storeKeys(json,jsonpath,var)"

# ---------------------------------------------------------------------
# 2) Remove the "text" column: shift columns Z..AE (26..31) left by one
#    into Y..AD (25..30) for every data row, then blank out column AE.
# ---------------------------------------------------------------------
for ($row = 1; $row -le 129; $row++) {
    for ($col = 26; $col -le 31; $col++) {
        $val = $ws.Cells.Item($row, $col).Value()
        $ws.Cells.Item($row, $col - 1).Value = $val
    }
    $ws.Cells.Item($row, 31).Value = ""
}

# ---------------------------------------------------------------------
# 3) Remove the "text" entry from the target/category list in column A:
#    shift A26:A31 up to A25:A30, then blank out A31.
# ---------------------------------------------------------------------
for ($row = 25; $row -le 30; $row++) {
    $val = $ws.Cells.Item($row + 1, 1).Value()
    $ws.Cells.Item($row, 1).Value = $val
}
$ws.Cells.Item(31, 1).Value = ""

# ---------------------------------------------------------------------
# 4) Update the defined names that describe the resized/shifted ranges.
# ---------------------------------------------------------------------
$wb.Names.Item("json").RefersTo = "='#system'!`$M`$2:`$M`$18"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AD`$2:`$AD`$27"
